$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50: fill in F, G, I, J (reusing existing shared strings)
$ws.Range("F50").Value = "Number of Troops"
$ws.Range("G50").Value = "March, 2013"
$ws.Range("I50").Value = "Afghanistan"
$ws.Range("J50").Value = "National"

# Row 53: fill in F, G, H, I, J (F, H, I are brand new strings)
$ws.Range("F53").Value = "People Affected, People Killed, Houses Affected"
$ws.Range("G53").Value = "March, 2013"
$ws.Range("H53").Value = "Government of Bolivia (Ministerio de Defensa)"
$ws.Range("I53").Value = "Bolivia"
$ws.Range("J53").Value = "Village / City"

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("H7").Select()
